$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "58.043.64"
Set-TextValue "E2" "  -0.99%  "
Set-TextValue "D3" "2.580.90"
Set-TextValue "E3" "  -1.72%  "
Set-TextValue "D5" "515.32"
Set-TextValue "E5" "  -1.76%  "
Set-TextValue "D6" "140.82"
Set-TextValue "E6" "  -1.63%  "
Set-TextValue "D7" "0.998"
Set-TextValue "E7" "  -0.09%  "
Set-TextValue "D8" "0.564"
Set-TextValue "E8" "  -1.04%  "
Set-TextValue "D9" "2.603.18"
Set-TextValue "E9" "  -1.24%  "
Set-TextValue "D10" "6.46"
Set-TextValue "E10" "  -2.08%  "
Set-TextValue "D11" "0.100"
Set-TextValue "E11" "  -2.16%  "
Set-TextValue "D12" "0.333"
Set-TextValue "E12" "  -0.54%  "
Set-TextValue "E13" "  -0.42%  "
Set-TextValue "D14" "3.041.68"
Set-TextValue "E14" "  -1.42%  "
Set-TextValue "D15" "58.024.00"
Set-TextValue "E15" "  -0.90%  "
Set-TextValue "D16" "20.17"
Set-TextValue "E16" "  -3.21%  "
Set-TextValue "D17" "0.0000132"
Set-TextValue "E17" "  -2.52%  "
Set-TextValue "D18" "2.559.74"
Set-TextValue "E18" "  -4.49%  "
Set-TextValue "D19" "334.76"
Set-TextValue "E19" "  -1.05%  "
Set-TextValue "D20" "4.28"
Set-TextValue "E20" "  -3.34%  "
Set-TextValue "D21" "10.10"
Set-TextValue "E21" "  -3.20%  "
Set-TextValue "D22" "6.39"
Set-TextValue "E22" "  +1.31%  "
Set-TextValue "D23" "0.998"
Set-TextValue "E23" "  -0.12%  "
Set-TextValue "D24" "65.23"
Set-TextValue "E24" "  -0.23%  "
Set-TextValue "E25" "  +0.53%  "
Set-TextValue "D26" "0.400"
Set-TextValue "E26" "  -3.82%  "
Set-TextValue "E27" "  -0.17%  "
Set-TextValue "D28" "2.711.87"
Set-TextValue "E28" "  -1.21%  "
Set-TextValue "D29" "6.97"
Set-TextValue "E29" "  -2.90%  "
Set-TextValue "D30" "0.998"
Set-TextValue "E30" "  -0.03%  "
Set-TextValue "D31" "0.0₃0723"
Set-TextValue "E31" "  -8.81%  "
Set-TextValue "D32" "6.04"
Set-TextValue "E32" "  -7.13%  "
Set-TextValue "D33" "1.56"
Set-TextValue "E33" "  -2.54%  "
Set-TextValue "D34" "18.65"
Set-TextValue "E34" "  -1.15%  "
Set-TextValue "D35" "149.34"
Set-TextValue "E35" "  -0.40%  "
Set-TextValue "D36" "3.95"
Set-TextValue "E36" "  -4.18%  "
Set-TextValue "D37" "1.12"
Set-TextValue "E37" "  -5.46%  "
Set-TextValue "B38" "Fetch.AI"
Set-TextValue "C38" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D38" "0.839"
Set-TextValue "E38" "  -1.42%  "
Set-TextValue "B39" "OKB"
Set-TextValue "C39" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D39" "36.01"
Set-TextValue "E39" "  -1.01%  "
Set-TextValue "B40" "SuiNetwork"
Set-TextValue "C40" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D40" "0.829"
Set-TextValue "E40" "  -7.77%  "
Set-TextValue "D41" "1.44"
Set-TextValue "E41" "  +0.08%  "
Set-TextValue "D42" "3.50"
Set-TextValue "E42" "  -3.58%  "
Set-TextValue "D43" "0.998"
Set-TextValue "E43" "  +0.06%  "
Set-TextValue "D44" "0.601"
Set-TextValue "E44" "  -0.28%  "
Set-TextValue "D45" "267.70"
Set-TextValue "E45" "  -0.92%  "
Set-TextValue "E46" "  +0.08%  "
Set-TextValue "D47" "0.0947"
Set-TextValue "E47" "  -2.58%  "
Set-TextValue "D48" "18.56"
Set-TextValue "E48" "  -2.89%  "
Set-TextValue "D49" "0.0517"
Set-TextValue "E49" "  -3.22%  "
Set-TextValue "D50" "1.954.66"
Set-TextValue "E50" "  -4.13%  "
Set-TextValue "D51" "4.57"
